$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (pushes current rows 5 and 6 down to 6 and 7)
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the latest weekly price data
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 44482
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100107
$ws.Cells.Item(5, 8).Value = "Otros"
$ws.Cells.Item(5, 9).Value = 100107002
$ws.Cells.Item(5, 10).Value = "Chirimoya"
$ws.Cells.Item(5, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 25000
$ws.Cells.Item(5, 15).Value = 26000
$ws.Cells.Item(5, 16).Value = 25500
$ws.Cells.Item(5, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(5, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(5, 19).Value = 2125
$ws.Cells.Item(5, 20).Value = 12
